# Commit: "unify the conception of DataNode, DataTable, Entity."
# The core semantic change is renaming the worksheet from the legacy
# "Property1" label to "DataNode" (aligning the sheet naming with the
# DataNode/DataTable/Entity concepts mentioned in the commit message).
# The remaining hunks in the upstream diff are incidental churn produced
# by re-saving the workbook in a newer Excel build (refreshed revision
# GUIDs, pruned-but-unused style pool entries, minor row autofit height
# drift, and the last-saved cell selection) - we reproduce the handful of
# those that are both meaningful and safely reproducible through the
# Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab: Property1 -> DataNode
$ws.Name = "DataNode"

# 2) The workbook was last saved with the cursor resting on H33 (instead
#    of the prior K17) - move the live selection to match.
$ws.Range("H33").Select() | Out-Null

# 3) A handful of rows pick up slightly different (auto-fit driven) row
#    heights in the re-saved workbook. Apply the same deltas.
$ws.Rows.Item(8).RowHeight = 92.1

$rows14 = 172,173,174,175,176,177,178,179,180,181,182,183,184,185,186,187,188,189,190
foreach ($r in $rows14) {
    $ws.Rows.Item($r).RowHeight = 14.1
}

# Rows 280-306 (even rows only) lose their explicit 14.25 override and
# fall back to the sheet's default row height - AutoFit clears the
# "custom height" flag instead of just writing 13.5 explicitly.
$rowsAuto = 280,282,284,286,288,290,292,294,296,298,300,302,304,306
foreach ($r in $rowsAuto) {
    $ws.Rows.Item($r).AutoFit()
}

# 4) The re-saved workbook's style table only keeps the default "Normal"
#    cell style - all the other (unused) built-in named cell styles that
#    Excel had accumulated over time are dropped. Deleting them does not
#    touch any cellXfs actually referenced by cells, so formatting is
#    unaffected.
for ($i = $wb.Styles.Count; $i -ge 2; $i--) {
    $style = $wb.Styles.Item($i)
    try {
        $style.Delete()
    } catch {
    }
}
